# Refresh the cryptocurrency price/volume table (rows 2-51) with the latest
# data pulled by the scheduled GitHub Actions job.
#
# Columns D (Price) and E (Volume(1h)) are stored as plain text in this sheet
# (values like "26.356.00" or "0.0₃0724" are not valid numbers, and even the
# ones that look numeric, e.g. "0.999" or "1.00", must stay text so trailing
# zeros / formatting survive). Any new value that Excel would otherwise parse
# as a real number is entered with a leading apostrophe, exactly like a user
# forcing text entry in the UI, so it keeps its original "General" look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.356.00'
$ws.Range("E2").Value = '  +0.38%  '

# Row 3
$ws.Range("D3").Value = '1.611.03'
$ws.Range("E3").Value = '  +1.36%  '

# Row 4
$ws.Range("D4").Value = "'" + '0.999'
$ws.Range("E4").Value = '  -0.17%  '

# Row 5
$ws.Range("D5").Value = "'" + '213.69'
$ws.Range("E5").Value = '  +0.65%  '

# Row 6
$ws.Range("D6").Value = "'" + '0.501'
$ws.Range("E6").Value = '  +0.02%  '

# Row 7
$ws.Range("D7").Value = "'" + '1.00'

# Row 8
$ws.Range("D8").Value = "'" + '0.246'
$ws.Range("E8").Value = '  -0.02%  '

# Row 9
$ws.Range("E9").Value = '  +0.20%  '

# Row 10
$ws.Range("D10").Value = "'" + '19.14'
$ws.Range("E10").Value = '  -1.08%  '

# Row 11
$ws.Range("D11").Value = "'" + '0.0855'

# Row 12
$ws.Range("D12").Value = '1.837.92'
$ws.Range("E12").Value = '  +1.39%  '

# Row 13
$ws.Range("D13").Value = '1.591.16'
$ws.Range("E13").Value = '  -1.55%  '

# Row 14
$ws.Range("E14").Value = '  +0.27%  '

# Row 15
$ws.Range("D15").Value = "'" + '0.512'
$ws.Range("E15").Value = '  -1.59%  '

# Row 16
$ws.Range("D16").Value = "'" + '64.60'
$ws.Range("E16").Value = '  +0.47%  '

# Row 17
$ws.Range("D17").Value = '26.374.63'
$ws.Range("E17").Value = '  +0.46%  '

# Row 18
$ws.Range("D18").Value = '0.0₃0724'
$ws.Range("E18").Value = '  -0.25%  '

# Row 19
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = "'" + '220.57'
$ws.Range("E19").Value = '  +3.69%  '

# Row 20
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = "'" + '7.53'
$ws.Range("E20").Value = '  +1.69%  '

# Row 21
$ws.Range("E21").Value = '  -0.18%  '

# Row 22
$ws.Range("D22").Value = "'" + '4.36'
$ws.Range("E22").Value = '  +1.70%  '

# Row 23
$ws.Range("D23").Value = "'" + '9.07'
$ws.Range("E23").Value = '  +0.70%  '

# Row 24
$ws.Range("E24").Value = '  -0.41%  '

# Row 25
$ws.Range("D25").Value = "'" + '144.96'
$ws.Range("E25").Value = '  +0.43%  '

# Row 27
$ws.Range("D27").Value = "'" + '7.02'
$ws.Range("E27").Value = '  -0.53%  '

# Row 28
$ws.Range("E28").Value = '  +1.18%  '

# Row 29
$ws.Range("D29").Value = "'" + '15.28'
$ws.Range("E29").Value = '  +0.66%  '

# Row 30
$ws.Range("E30").Value = '  -0.24%  '

# Row 31
$ws.Range("E31").Value = '  +0.12%  '

# Row 32
$ws.Range("D32").Value = "'" + '3.21'
$ws.Range("E32").Value = '  +0.58%  '

# Row 33
$ws.Range("D33").Value = '1.443.16'
$ws.Range("E33").Value = '  +8.12%  '

# Row 34
$ws.Range("E34").Value = '  +1.41%  '

# Row 35
$ws.Range("E35").Value = '  -0.83%  '

# Row 36
$ws.Range("E36").Value = '  +0.08%  '

# Row 37
$ws.Range("E37").Value = '  -4.96%  '

# Row 38
$ws.Range("E38").Value = '  -0.17%  '

# Row 39
$ws.Range("D39").Value = "'" + '0.835'
$ws.Range("E39").Value = '  +2.15%  '

# Row 40
$ws.Range("E40").Value = '  +1.86%  '

# Row 41
$ws.Range("D41").Value = "'" + '1.00'
$ws.Range("E41").Value = '  -0.14%  '

# Row 42
$ws.Range("E42").Value = '  +1.94%  '

# Row 43
$ws.Range("D43").Value = '1.749.42'
$ws.Range("E43").Value = '  +1.29%  '

# Row 44
$ws.Range("D44").Value = "'" + '0.761'
$ws.Range("E44").Value = '  -0.11%  '

# Row 45
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = "'" + '61.73'
$ws.Range("E45").Value = '  -0.23%  '

# Row 46
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = "'" + '0.909'
$ws.Range("E46").Value = '  -11.33%  '

# Row 47
$ws.Range("D47").Value = "'" + '87.95'
$ws.Range("E47").Value = '  +2.67%  '

# Row 48
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = "'" + '1.49'
$ws.Range("E48").Value = '  +0.41%  '

# Row 49
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = "'" + '0.0502'
$ws.Range("E49").Value = '  -0.07%  '

# Row 50
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = "'" + '0.0960'
$ws.Range("E50").Value = '  -1.35%  '

# Row 51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = "'" + '7.47'
$ws.Range("E51").Value = '  +1.41%  '
